$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8949.210999999999
$ws.Range("I74").Value = 5670.3335
$ws.Range("K74").Value = 5670.3335
$ws.Range("M74").Value = -4734.3335

$ws.Range("H77").Value = 8949.210999999999
$ws.Range("I77").Value = 5670.3335
$ws.Range("K77").Value = 28351.6675
$ws.Range("M77").Value = -23671.6675

$ws.Range("H98").Value = 1894.8649
$ws.Range("I98").Value = 1177.9
$ws.Range("K98").Value = 1177.9
$ws.Range("M98").Value = 320.0999999999999

$ws.Range("H100").Value = 12562.167
$ws.Range("I100").Value = 872.5
$ws.Range("J100").Value = 14900.1
$ws.Range("K100").Value = 872.5
$ws.Range("L100").Value = 14900.1
$ws.Range("M100").Value = -331.5
$ws.Range("N100").Value = -15982.1

$ws.Range("H116").Value = 16962.777
$ws.Range("J116").Value = 25111.6
$ws.Range("L116").Value = 25111.6
$ws.Range("N116").Value = -31995.6

$ws.Range("H122").Value = 1894.8649
$ws.Range("I122").Value = 1177.9
$ws.Range("K122").Value = 3533.7
$ws.Range("M122").Value = -1083.7

$ws.Range("H125").Value = 2659.6365
$ws.Range("I125").Value = 2086
$ws.Range("J125").Value = 3137.6667
$ws.Range("K125").Value = 18774
$ws.Range("L125").Value = 28239.0003
$ws.Range("M125").Value = -16314
$ws.Range("N125").Value = -33159.0003

$ws.Range("H135").Value = 1668338.4
$ws.Range("I135").Value = 2501398
$ws.Range("K135").Value = 22512582
$ws.Range("M135").Value = -22510047

$ws.Range("H137").Value = 1435724.9
$ws.Range("I137").Value = 1113349.5
$ws.Range("J137").Value = 2016000.6
$ws.Range("K137").Value = 3340048.5
$ws.Range("L137").Value = 6048001.800000001
$ws.Range("M137").Value = -3337498.5
$ws.Range("N137").Value = -6053101.800000001

$ws.Range("H138").Value = 7629.675
$ws.Range("J138").Value = 12051.087
$ws.Range("L138").Value = 36153.261
$ws.Range("N138").Value = -46433.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3291.5
$ws.Range("J45").Value = 4766.6665
$ws.Range("L45").Value = 4766.6665
$ws.Range("N45").Value = -5520.6665

$ws.Range("H122").Value = 3933.8333
$ws.Range("I122").Value = 2158.4666
$ws.Range("K122").Value = 6475.399800000001
$ws.Range("M122").Value = -4025.399800000001

$ws.Range("H132").Value = 3924.5208
$ws.Range("I132").Value = 2927.0789
$ws.Range("K132").Value = 8781.236699999999
$ws.Range("M132").Value = -6251.236699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 36943.5
$ws.Range("I134").Value = 2502.9473
$ws.Range("K134").Value = 7508.841899999999
$ws.Range("M134").Value = -4973.841899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1481
$ws.Range("I16").Value = 2072.6
$ws.Range("K16").Value = 2072.6
$ws.Range("M16").Value = -1785.6

$ws.Range("H58").Value = 5839.524
$ws.Range("I58").Value = 4793.5386
$ws.Range("K58").Value = 4793.5386
$ws.Range("M58").Value = -4590.5386

$ws.Range("H113").Value = 1481
$ws.Range("I113").Value = 2072.6
$ws.Range("K113").Value = 2072.6
$ws.Range("M113").Value = 97.40000000000009

$ws.Range("H136").Value = 5839.524
$ws.Range("I136").Value = 4793.5386
$ws.Range("K136").Value = 14380.6158
$ws.Range("M136").Value = -11830.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 600398.4
$ws.Range("I5").Value = 53677.734
$ws.Range("J5").Value = 3334001.8
$ws.Range("K5").Value = 161033.202
$ws.Range("L5").Value = 10002005.4
$ws.Range("M5").Value = -160921.202
$ws.Range("N5").Value = -10002229.4

$ws.Range("H55").Value = 47405.883
$ws.Range("J55").Value = 49993.75
$ws.Range("L55").Value = 149981.25
$ws.Range("N55").Value = -150335.25

$ws.Range("H68").Value = 1429284.2
$ws.Range("I68").Value = 1429013
$ws.Range("J68").Value = 1429555.6
$ws.Range("K68").Value = 4287039
$ws.Range("L68").Value = 4288666.800000001
$ws.Range("M68").Value = -4286228
$ws.Range("N68").Value = -4290288.800000001

$ws.Range("H71").Value = 1429284.2
$ws.Range("I71").Value = 1429013
$ws.Range("J71").Value = 1429555.6
$ws.Range("K71").Value = 12861117
$ws.Range("L71").Value = 12866000.4
$ws.Range("M71").Value = -12857061
$ws.Range("N71").Value = -12874112.4

$ws.Range("H97").Value = 1119.8
$ws.Range("J97").Value = 899
$ws.Range("L97").Value = 2697
$ws.Range("N97").Value = -3689

$ws.Range("H135").Value = 600398.4
$ws.Range("I135").Value = 53677.734
$ws.Range("J135").Value = 3334001.8
$ws.Range("K135").Value = 483099.606
$ws.Range("L135").Value = 30006016.2
$ws.Range("M135").Value = -480564.606
$ws.Range("N135").Value = -30011086.2

$ws.Range("H139").Value = 4567.533
$ws.Range("I139").Value = 3550.4783
$ws.Range("K139").Value = 10651.4349
$ws.Range("M139").Value = -5511.4349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 564364.75
$ws.Range("I113").Value = 1431414.6
$ws.Range("J113").Value = 12605.728
$ws.Range("K113").Value = 1431414.6
$ws.Range("L113").Value = 12605.728
$ws.Range("M113").Value = -1429244.6
$ws.Range("N113").Value = -16945.728

$ws.Range("H122").Value = 3577.2104
$ws.Range("J122").Value = 5216.3335
$ws.Range("L122").Value = 15649.0005
$ws.Range("N122").Value = -20549.0005

$ws.Range("H126").Value = 4398.1
$ws.Range("I126").Value = 3381
$ws.Range("K126").Value = 10143
$ws.Range("M126").Value = -7673

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 481605.66
$ws.Range("I7").Value = 6424.8
$ws.Range("J7").Value = 913588.25
$ws.Range("K7").Value = 6424.8
$ws.Range("L7").Value = 913588.25
$ws.Range("M7").Value = -6312.8
$ws.Range("N7").Value = -913812.25

$ws.Range("H40").Value = 1353576.2
$ws.Range("I40").Value = 1787531.1
$ws.Range("J40").Value = 3494.4443
$ws.Range("K40").Value = 1787531.1
$ws.Range("L40").Value = 3494.4443
$ws.Range("M40").Value = -1787395.1
$ws.Range("N40").Value = -3766.4443

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""

$ws.Range("H68").Value = 3379.1667
$ws.Range("J68").Value = 3675.5
$ws.Range("L68").Value = 3675.5
$ws.Range("N68").Value = -5173.5

$ws.Range("H71").Value = 3379.1667
$ws.Range("J71").Value = 3675.5
$ws.Range("L71").Value = 18377.5
$ws.Range("N71").Value = -25865.5

$ws.Range("H93").Value = 1970.4706
$ws.Range("I93").Value = 1137.5834
$ws.Range("J93").Value = 3969.4
$ws.Range("K93").Value = 1137.5834
$ws.Range("L93").Value = 3969.4
$ws.Range("M93").Value = 110.4166
$ws.Range("N93").Value = -6465.4

$ws.Range("H126").Value = 481605.66
$ws.Range("I126").Value = 6424.8
$ws.Range("J126").Value = 913588.25
$ws.Range("K126").Value = 19274.4
$ws.Range("L126").Value = 2740764.75
$ws.Range("M126").Value = -16804.4
$ws.Range("N126").Value = -2745704.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""

$ws.Range("H12").Value = 876.75
$ws.Range("J12").Value = 876.75
$ws.Range("L12").Value = 876.75
$ws.Range("N12").Value = -1160.75

$ws.Range("H81").Value = 1932
$ws.Range("J81").Value = 2883.5
$ws.Range("L81").Value = 5767
$ws.Range("N81").Value = -7889

$ws.Range("H84").Value = 1932
$ws.Range("J84").Value = 2883.5
$ws.Range("L84").Value = 28835
$ws.Range("N84").Value = -39443

$ws.Range("H100").Value = 690.75
$ws.Range("I100").Value = 695.5263
$ws.Range("K100").Value = 1391.0526
$ws.Range("M100").Value = -850.0526
